$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a brand-new "Problems with meta-layers" Heading2 section
#    (heading + 2 body paragraphs) immediately before the existing
#    "Entering layer data" Heading2 section.
# ------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("Entering layer data", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$headingPara = $find.Paragraphs(1)
$insertPoint = $headingPara.Range.Duplicate
$insertPoint.Collapse(1)

$newText = "Problems with meta-layers`r" `
  + "`r" `
  + "We need to think carefully when we should allow meta-layers and when we shouldn" + [char]0x2019 + "t. For example meta-layers under biomass events might be misleading if only applied for specific years when it should be all years. `r" `
  + "Maybe have meta-layers used instead of layers as default, but with specific safeguards or cases where you can" + [char]0x2019 + "t.`r"

$insertPoint.InsertBefore($newText)

# ------------------------------------------------------------------
# 2. Fix up paragraph styles: the inserted paragraphs all inherited
#    the Heading2 style from the split point, so restore the body
#    paragraphs to Normal and keep only the true heading as Heading2.
# ------------------------------------------------------------------
$rngHeading = $d.Content
$rngHeading.Find.Execute("Problems with meta-layers", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngHeading.Paragraphs(1).Range.Style = "Heading 2"

$rngP1 = $d.Content
$rngP1.Find.Execute("We need to think carefully", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngP1.Paragraphs(1).Range.Style = "Normal"

$rngP2 = $d.Content
$rngP2.Find.Execute("Maybe have meta-layers used", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngP2.Paragraphs(1).Range.Style = "Normal"

# ------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark: it used to sit at the end of the
#    "... It works." paragraph; it now belongs at the end of the new
#    "Maybe have meta-layers ..." paragraph (the author's last edit).
# ------------------------------------------------------------------
$rngBm = $d.Content
$rngBm.Find.Execute("Maybe have meta-layers used instead of layers as default, but with specific safeguards or cases where you can" + [char]0x2019 + "t.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngBm.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rngBm) | Out-Null

Write-Output "done"
